$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 22.81767259378528
$ws.Cells.Item(2, 5).Value = 20.67592430114746
$ws.Cells.Item(2, 6).Value = 24.17180914708917
$ws.Cells.Item(2, 7).Value = 20.65289510857373
$ws.Cells.Item(2, 8).Value = 167899010
$ws.Cells.Item(2, 9).Value = "CGNX"

$ws.Cells.Item(3, 4).Value = 22.35593356154076
$ws.Cells.Item(3, 5).Value = 20.87998962402344
$ws.Cells.Item(3, 6).Value = 22.48046651271329
$ws.Cells.Item(3, 7).Value = 20.36802141124588
$ws.Cells.Item(3, 8).Value = 167899010
$ws.Cells.Item(3, 9).Value = "CGNX"

$ws.Cells.Item(4, 4).Value = 15.89887802760023
$ws.Cells.Item(4, 5).Value = 17.37784194946289
$ws.Cells.Item(4, 6).Value = 17.46103389045067
$ws.Cells.Item(4, 7).Value = 15.67241068534147
$ws.Cells.Item(4, 8).Value = 167899010
$ws.Cells.Item(4, 9).Value = "CGNX"

$ws.Cells.Item(5, 4).Value = 15.29971221972721
$ws.Cells.Item(5, 5).Value = 14.93388938903809
$ws.Cells.Item(5, 6).Value = 15.50809242495885
$ws.Cells.Item(5, 7).Value = 12.97048822341249
$ws.Cells.Item(5, 8).Value = 167899010
$ws.Cells.Item(5, 9).Value = "CGNX"

$ws.Cells.Item(6, 4).Value = 18.05545739094218
$ws.Cells.Item(6, 5).Value = 16.48279571533203
$ws.Cells.Item(6, 6).Value = 18.32916559024523
$ws.Cells.Item(6, 7).Value = 16.3992913989559
$ws.Cells.Item(6, 8).Value = 167899010
$ws.Cells.Item(6, 9).Value = "CGNX"

$ws.Cells.Item(7, 4).Value = 19.9551433668408
$ws.Cells.Item(7, 5).Value = 20.99147224426269
$ws.Cells.Item(7, 6).Value = 21.07977047620426
$ws.Cells.Item(7, 7).Value = 19.48577536611787
$ws.Cells.Item(7, 8).Value = 167899010
$ws.Cells.Item(7, 9).Value = "CGNX"

$ws.Cells.Item(8, 4).Value = 24.6066801136449
$ws.Cells.Item(8, 5).Value = 24.01559829711914
$ws.Cells.Item(8, 6).Value = 25.49097613610129
$ws.Cells.Item(8, 7).Value = 23.68049812291513
$ws.Cells.Item(8, 8).Value = 167899010
$ws.Cells.Item(8, 9).Value = "CGNX"

$ws.Cells.Item(9, 4).Value = 29.55357674002671
$ws.Cells.Item(9, 5).Value = 31.48280715942383
$ws.Cells.Item(9, 6).Value = 31.66920697229829
$ws.Cells.Item(9, 7).Value = 29.059620257898
$ws.Cells.Item(9, 8).Value = 167899010
$ws.Cells.Item(9, 9).Value = "CGNX"

$ws.Cells.Item(10, 4).Value = 39.74181048536277
$ws.Cells.Item(10, 5).Value = 39.80711364746094
$ws.Cells.Item(10, 6).Value = 40.52545198929624
$ws.Cells.Item(10, 7).Value = 37.06903569143721
$ws.Cells.Item(10, 8).Value = 167899010
$ws.Cells.Item(10, 9).Value = "CGNX"

$ws.Cells.Item(11, 4).Value = 39.95561458374878
$ws.Cells.Item(11, 5).Value = 44.38163757324219
$ws.Cells.Item(11, 6).Value = 45.08195798213794
$ws.Cells.Item(11, 7).Value = 39.43270739611498
$ws.Cells.Item(11, 8).Value = 167899010
$ws.Cells.Item(11, 9).Value = "CGNX"

$ws.Cells.Item(12, 4).Value = 51.72128323898426
$ws.Cells.Item(12, 5).Value = 57.54337310791016
$ws.Cells.Item(12, 6).Value = 60.69272111453528
$ws.Cells.Item(12, 7).Value = 51.63717585293529
$ws.Cells.Item(12, 8).Value = 167899010
$ws.Cells.Item(12, 9).Value = "CGNX"

$ws.Cells.Item(13, 4).Value = 57.333231726735
$ws.Cells.Item(13, 5).Value = 58.3244743347168
$ws.Cells.Item(13, 6).Value = 65.45956827558983
$ws.Cells.Item(13, 7).Value = 56.92176930023987
$ws.Cells.Item(13, 8).Value = 167899010
$ws.Cells.Item(13, 9).Value = "CGNX"

$ws.Cells.Item(14, 4).Value = 48.42455431848412
$ws.Cells.Item(14, 5).Value = 43.28634643554688
$ws.Cells.Item(14, 6).Value = 48.6398157209762
$ws.Cells.Item(14, 7).Value = 42.21003585283116
$ws.Cells.Item(14, 8).Value = 167899010
$ws.Cells.Item(14, 9).Value = "CGNX"

$ws.Cells.Item(15, 4).Value = 41.584293711384
$ws.Cells.Item(15, 5).Value = 49.44399642944336
$ws.Cells.Item(15, 6).Value = 50.98970827940524
$ws.Cells.Item(15, 7).Value = 40.9660096861165
$ws.Cells.Item(15, 8).Value = 167899010
$ws.Cells.Item(15, 9).Value = "CGNX"

$ws.Cells.Item(16, 4).Value = 52.56191488185687
$ws.Cells.Item(16, 5).Value = 40.16682815551758
$ws.Cells.Item(16, 6).Value = 53.5839001447547
$ws.Cells.Item(16, 7).Value = 33.90365300656529
$ws.Cells.Item(16, 8).Value = 167899010
$ws.Cells.Item(16, 9).Value = "CGNX"

$ws.Cells.Item(17, 4).Value = 35.40921017724591
$ws.Cells.Item(17, 5).Value = 42.71259307861328
$ws.Cells.Item(17, 6).Value = 43.00360324369439
$ws.Cells.Item(17, 7).Value = 33.41908238295358
$ws.Cells.Item(17, 8).Value = 167899010
$ws.Cells.Item(17, 9).Value = "CGNX"

$ws.Cells.Item(18, 4).Value = 48.63407925333165
$ws.Cells.Item(18, 5).Value = 47.38440322875977
$ws.Cells.Item(18, 6).Value = 53.84890338448178
$ws.Cells.Item(18, 7).Value = 42.25414749444127
$ws.Cells.Item(18, 8).Value = 167899010
$ws.Cells.Item(18, 9).Value = "CGNX"

$ws.Cells.Item(19, 4).Value = 46.2137752278483
$ws.Cells.Item(19, 5).Value = 41.39768218994141
$ws.Cells.Item(19, 6).Value = 46.84400409228045
$ws.Cells.Item(19, 7).Value = 37.6068933907163
$ws.Cells.Item(19, 8).Value = 167899010
$ws.Cells.Item(19, 9).Value = "CGNX"

$ws.Cells.Item(20, 4).Value = 46.65276375714611
$ws.Cells.Item(20, 5).Value = 48.48911666870117
$ws.Cells.Item(20, 6).Value = 50.9846700577199
$ws.Cells.Item(20, 7).Value = 43.4791678933651
$ws.Cells.Item(20, 8).Value = 167899010
$ws.Cells.Item(20, 9).Value = "CGNX"

$ws.Cells.Item(21, 4).Value = 53.20790361109589
$ws.Cells.Item(21, 5).Value = 48.05115127563477
$ws.Cells.Item(21, 6).Value = 55.75328603207248
$ws.Cells.Item(21, 7).Value = 47.68348560522881
$ws.Cells.Item(21, 8).Value = 167899010
$ws.Cells.Item(21, 9).Value = "CGNX"

$ws.Cells.Item(22, 4).Value = 38.35785542324881
$ws.Cells.Item(22, 5).Value = 52.13799285888672
$ws.Cells.Item(22, 6).Value = 55.28099444327352
$ws.Cells.Item(22, 7).Value = 37.82930459370862
$ws.Cells.Item(22, 8).Value = 167899010
$ws.Cells.Item(22, 9).Value = "CGNX"

$ws.Cells.Item(23, 4).Value = 56.33830004501073
$ws.Cells.Item(23, 5).Value = 63.17863845825195
$ws.Cells.Item(23, 6).Value = 63.25421681233991
$ws.Cells.Item(23, 7).Value = 55.36515575273921
$ws.Cells.Item(23, 8).Value = 167899010
$ws.Cells.Item(23, 9).Value = "CGNX"

$ws.Cells.Item(24, 4).Value = 62.24739702201942
$ws.Cells.Item(24, 5).Value = 62.3135871887207
$ws.Cells.Item(24, 6).Value = 67.07930083377447
$ws.Cells.Item(24, 7).Value = 59.66596609229513
$ws.Cells.Item(24, 8).Value = 167899010
$ws.Cells.Item(24, 9).Value = "CGNX"

$ws.Cells.Item(25, 4).Value = 78.71029603509299
$ws.Cells.Item(25, 5).Value = 79.81820678710938
$ws.Cells.Item(25, 6).Value = 85.66876032354666
$ws.Cells.Item(25, 7).Value = 76.77631251740023
$ws.Cells.Item(25, 8).Value = 167899010
$ws.Cells.Item(25, 9).Value = "CGNX"

$ws.Cells.Item(26, 4).Value = 81.50931238251138
$ws.Cells.Item(26, 5).Value = 83.75590515136719
$ws.Cells.Item(26, 6).Value = 86.42068678018498
$ws.Cells.Item(26, 7).Value = 80.0213139088283
$ws.Cells.Item(26, 8).Value = 167899010
$ws.Cells.Item(26, 9).Value = "CGNX"

$ws.Cells.Item(27, 4).Value = 81.36927446076839
$ws.Cells.Item(27, 5).Value = 87.99756622314453
$ws.Cells.Item(27, 6).Value = 88.47448932589549
$ws.Cells.Item(27, 7).Value = 78.71212403391024
$ws.Cells.Item(27, 8).Value = 167899010
$ws.Cells.Item(27, 9).Value = "CGNX"

$ws.Cells.Item(28, 4).Value = 78.29114299876429
$ws.Cells.Item(28, 5).Value = 85.31377410888672
$ws.Cells.Item(28, 6).Value = 86.20987149603181
$ws.Cells.Item(28, 7).Value = 77.10285367332946
$ws.Cells.Item(28, 8).Value = 167899010
$ws.Cells.Item(28, 9).Value = "CGNX"

$ws.Cells.Item(29, 4).Value = 75.95571055042126
$ws.Cells.Item(29, 5).Value = 64.78460693359375
$ws.Cells.Item(29, 6).Value = 76.42361348803976
$ws.Cells.Item(29, 7).Value = 60.08611475858883
$ws.Cells.Item(29, 8).Value = 167899010
$ws.Cells.Item(29, 9).Value = "CGNX"

$ws.Cells.Item(30, 4).Value = 75.13080142263452
$ws.Cells.Item(30, 5).Value = 65.98825836181641
$ws.Cells.Item(30, 6).Value = 77.52132394477016
$ws.Cells.Item(30, 7).Value = 65.83214628178642
$ws.Cells.Item(30, 8).Value = 167899010
$ws.Cells.Item(30, 9).Value = "CGNX"

$ws.Cells.Item(31, 4).Value = 41.64230086655685
$ws.Cells.Item(31, 5).Value = 49.81052398681641
$ws.Cells.Item(31, 6).Value = 50.10364127453335
$ws.Cells.Item(31, 7).Value = 41.28078966927892
$ws.Cells.Item(31, 8).Value = 167899010
$ws.Cells.Item(31, 9).Value = "CGNX"

$ws.Cells.Item(32, 4).Value = 40.84834678296749
$ws.Cells.Item(32, 5).Value = 45.23159408569336
$ws.Cells.Item(32, 6).Value = 46.37632813993508
$ws.Cells.Item(32, 7).Value = 39.34160447294238
$ws.Cells.Item(32, 8).Value = 167899010
$ws.Cells.Item(32, 9).Value = "CGNX"

$ws.Cells.Item(33, 4).Value = 46.59788476978674
$ws.Cells.Item(33, 5).Value = 53.63263702392578
$ws.Cells.Item(33, 6).Value = 53.69142085994539
$ws.Cells.Item(33, 7).Value = 46.03941590243409
$ws.Cells.Item(33, 8).Value = 167899010
$ws.Cells.Item(33, 9).Value = "CGNX"

$ws.Cells.Item(34, 4).Value = 48.39275988647957
$ws.Cells.Item(34, 5).Value = 46.79340362548828
$ws.Cells.Item(34, 6).Value = 48.94223296355696
$ws.Cells.Item(34, 7).Value = 45.97900997260158
$ws.Cells.Item(34, 8).Value = 167899010
$ws.Cells.Item(34, 9).Value = "CGNX"

$ws.Cells.Item(35, 4).Value = 54.83457106546705
$ws.Cells.Item(35, 5).Value = 53.66536712646485
$ws.Cells.Item(35, 6).Value = 58.46990058720208
$ws.Cells.Item(35, 7).Value = 52.62389588429162
$ws.Cells.Item(35, 8).Value = 167899010
$ws.Cells.Item(35, 9).Value = "CGNX"

$ws.Cells.Item(36, 4).Value = 41.66170302498773
$ws.Cells.Item(36, 5).Value = 35.41343307495117
$ws.Cells.Item(36, 6).Value = 42.26193127680494
$ws.Cells.Item(36, 7).Value = 33.77018493033809
$ws.Cells.Item(36, 8).Value = 167899010
$ws.Cells.Item(36, 9).Value = "CGNX"

$ws.Cells.Item(37, 4).Value = 40.66092820202534
$ws.Cells.Item(37, 5).Value = 35.63253784179688
$ws.Cells.Item(37, 6).Value = 40.75952254273976
$ws.Cells.Item(37, 7).Value = 35.49450426034349
$ws.Cells.Item(37, 8).Value = 167899010
$ws.Cells.Item(37, 9).Value = "CGNX"

$ws.Cells.Item(38, 4).Value = 41.97482885159913
$ws.Cells.Item(38, 5).Value = 41.03634643554688
$ws.Cells.Item(38, 6).Value = 42.08349198385204
$ws.Cells.Item(38, 7).Value = 37.83562837105549
$ws.Cells.Item(38, 8).Value = 167899010
$ws.Cells.Item(38, 9).Value = "CGNX"

$ws.Cells.Item(39, 4).Value = 46.5042115042877
$ws.Cells.Item(39, 5).Value = 49.09657287597656
$ws.Cells.Item(39, 6).Value = 52.56954908286514
$ws.Cells.Item(39, 7).Value = 45.72254406462249
$ws.Cells.Item(39, 8).Value = 167899010
$ws.Cells.Item(39, 9).Value = "CGNX"

$ws.Cells.Item(40, 4).Value = 39.85356052891638
$ws.Cells.Item(40, 5).Value = 39.88330078125
$ws.Cells.Item(40, 6).Value = 43.78934483033753
$ws.Cells.Item(40, 7).Value = 37.68243134482139
$ws.Cells.Item(40, 8).Value = 167899010
$ws.Cells.Item(40, 9).Value = "CGNX"

$ws.Cells.Item(41, 4).Value = 35.68924167842574
$ws.Cells.Item(41, 5).Value = 39.63264083862305
$ws.Cells.Item(41, 6).Value = 41.15238744976228
$ws.Cells.Item(41, 7).Value = 35.29192047031056
$ws.Cells.Item(41, 8).Value = 167899010
$ws.Cells.Item(41, 9).Value = "CGNX"

$ws.Cells.Item(42, 4).Value = 29.70170834603361
$ws.Cells.Item(42, 5).Value = 27.18258857727051
$ws.Cells.Item(42, 6).Value = 30.65758060697582
$ws.Cells.Item(42, 7).Value = 22.5725019172909
$ws.Cells.Item(42, 8).Value = 167899010
$ws.Cells.Item(42, 9).Value = "CGNX"

$ws.Cells.Item(43, 4).Value = 31.47290653311684
$ws.Cells.Item(43, 5).Value = 40.69617462158203
$ws.Cells.Item(43, 6).Value = 42.06369277555218
$ws.Cells.Item(43, 7).Value = 31.44296017136102
$ws.Cells.Item(43, 8).Value = 167899010
$ws.Cells.Item(43, 9).Value = "CGNX"
